$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 6 (shifting the existing rows 6-10 down to 7-11),
# copying the formatting of the row above (row 5) so borders/alignment match.
$ws.Rows.Item(6).Insert()
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)

# Populate the new "Lasso Regression+normalization+ lag1" result row.
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Lasso Regression+normalization+ lag1"
$ws.Range("C6").Value = 77.074776640014207

# Renumber the Id column for the rows that shifted down.
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9

# Move the active selection as recorded after the edit.
$ws.Range("B6").Select()
